$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers must be forced to Text format
# so Excel does not auto-convert them (matching the source data which stores
# prices/percentages as text strings).

$ws.Range("D2").Value = "89.455.14"
$ws.Range("E2").Value = "  -1.55%  "
$ws.Range("D3").Value = "3.138.70"
$ws.Range("E3").Value = "  -1.79%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.49"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "636.02"
$ws.Range("E6").Value = "  +3.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.396"
$ws.Range("E7").Value = "  +1.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.783"
$ws.Range("E8").Value = "  +13.58%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "3.136.54"
$ws.Range("E10").Value = "  -1.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.561"
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.34"
$ws.Range("E14").Value = "  +2.42%  "
$ws.Range("D15").Value = "89.257.95"
$ws.Range("E15").Value = "  -1.43%  "
$ws.Range("D16").Value = "3.713.04"
$ws.Range("E16").Value = "  -3.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "32.28"
$ws.Range("E17").Value = "  -2.25%  "
$ws.Range("D18").Value = "3.138.79"
$ws.Range("E18").Value = "  -3.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.42"
$ws.Range("E19").Value = "  +5.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000227"
$ws.Range("E20").Value = "  +19.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.26"
$ws.Range("E21").Value = "  -1.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "425.67"
$ws.Range("E22").Value = "  -2.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.43"
$ws.Range("E23").Value = "  -1.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.92"
$ws.Range("E24").Value = "  -3.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.45"
$ws.Range("E25").Value = "  +5.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "82.28"
$ws.Range("E26").Value = "  +9.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.51"
$ws.Range("E27").Value = "  -2.51%  "
$ws.Range("D28").Value = "3.299.68"
$ws.Range("E28").Value = "  -4.21%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("B30").Value = "Cronos"
$ws.Range("C30").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.159"
$ws.Range("E30").Value = "  -6.06%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.02"
$ws.Range("E32").Value = "  -3.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.18"
$ws.Range("E33").Value = "  -3.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "505.04"
$ws.Range("E34").Value = "  -5.96%  "
$ws.Range("E35").Value = "  +16.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.97"
$ws.Range("E36").Value = "  +0.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.28"
$ws.Range("E37").Value = "  +3.14%  "
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.24"
$ws.Range("E39").Value = "  +0.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.27"
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  -3.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.366"
$ws.Range("E44").Value = "  -2.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "145.56"
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.132"
$ws.Range("E46").Value = "  +7.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.72"
$ws.Range("E47").Value = "  -2.40%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "164.66"
$ws.Range("E48").Value = "  -5.42%  "
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0652"
$ws.Range("E49").Value = "  +11.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.723"
$ws.Range("E50").Value = "  +2.50%  "
$ws.Range("E51").Value = "  -0.09%  "
